# 25Tto26TMap.xlsx: remove the PrEP_CT.T and PrEP_CT.TestResult.T rows
# (commit: "removed PrEP_CT.R, PrEP_CT.T_1, PrEP_CT.T, PrEP_CT.TestResult.T")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

# Delete the higher-numbered row first so the lower row index is unaffected.
$ws.Rows.Item(48).Delete()   # "PrEP_CT.TestResult.T"
$ws.Rows.Item(24).Delete()   # "PrEP_CT.T"

# Re-establish the AutoFilter over the now-smaller used range (A1:N55) and
# keep the workbook-level _FilterDatabase name in sync with it.
$ws.AutoFilterMode = $false
$ws.Range("A1:N55").AutoFilter()
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Map!`$A`$1:`$N`$55"

# Restore the view: frozen header row with the body scrolled down to A27,
# and the selected cell sitting just past the last data row.
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("A58").Select()
